$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "51.937.36"
Set-TextValue 2 5 "  +2.07%  "

# Row 3
Set-TextValue 3 4 "3.007.73"
Set-TextValue 3 5 "  +3.58%  "

# Row 4
Set-TextValue 4 4 "0.998"
Set-TextValue 4 5 "  -0.03%  "

# Row 5
Set-TextValue 5 4 "385.27"
Set-TextValue 5 5 "  +3.59%  "

# Row 6
Set-TextValue 6 4 "105.22"
Set-TextValue 6 5 "  +3.55%  "

# Row 7
Set-TextValue 7 4 "0.549"
Set-TextValue 7 5 "  +1.47%  "

# Row 8
Set-TextValue 8 5 "  +0.02%  "

# Row 9
Set-TextValue 9 4 "0.601"
Set-TextValue 9 5 "  +2.99%  "

# Row 10
Set-TextValue 10 4 "37.60"
Set-TextValue 10 5 "  +2.31%  "

# Row 11
Set-TextValue 11 5 "  +0.39%  "

# Row 12
Set-TextValue 12 4 "0.0851"
Set-TextValue 12 5 "  +2.38%  "

# Row 13
Set-TextValue 13 4 "3.471.41"
Set-TextValue 13 5 "  +3.36%  "

# Row 14
Set-TextValue 14 4 "18.55"
Set-TextValue 14 5 "  +1.88%  "

# Row 15
Set-TextValue 15 4 "7.64"
Set-TextValue 15 5 "  +4.16%  "

# Row 16
Set-TextValue 16 2 "WrappedEther"
Set-TextValue 16 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 16 4 "3.004.37"
Set-TextValue 16 5 "  +3.49%  "

# Row 17
Set-TextValue 17 2 "Polygon"
Set-TextValue 17 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue 17 4 "1.02"
Set-TextValue 17 5 "  +11.06%  "

# Row 18
Set-TextValue 18 4 "51.785.34"
Set-TextValue 18 5 "  +1.92%  "

# Row 19
Set-TextValue 19 4 "3.32"
Set-TextValue 19 5 "  +2.70%  "

# Row 20
Set-TextValue 20 4 "7.49"
Set-TextValue 20 5 "  +4.38%  "

# Row 21
Set-TextValue 21 4 "13.04"
Set-TextValue 21 5 "  +1.58%  "

# Row 22
Set-TextValue 22 4 "0.0₃0970"
Set-TextValue 22 5 "  +3.43%  "

# Row 23
Set-TextValue 23 4 "69.30"

# Row 24
Set-TextValue 24 4 "264.55"
Set-TextValue 24 5 "  +2.17%  "

# Row 25
Set-TextValue 25 5 "  +9.85%  "

# Row 26
Set-TextValue 26 4 "8.46"
Set-TextValue 26 5 "  +20.24%  "

# Row 27
Set-TextValue 27 4 "7.67"
Set-TextValue 27 5 "  +23.38%  "

# Row 28
Set-TextValue 28 5 "  +1.97%  "

# Row 29
Set-TextValue 29 5 "  +13.70%  "

# Row 30
Set-TextValue 30 4 "26.22"
Set-TextValue 30 5 "  +2.71%  "

# Row 31
Set-TextValue 31 5 "  +0.02%  "

# Row 32
Set-TextValue 32 4 "9.95"
Set-TextValue 32 5 "  +1.28%  "

# Row 33
Set-TextValue 33 4 "35.28"
Set-TextValue 33 5 "  +3.93%  "

# Row 34
Set-TextValue 34 4 "51.18"
Set-TextValue 34 5 "  +0.03%  "

# Row 35
Set-TextValue 35 5 "  -1.67%  "

# Row 36
Set-TextValue 36 4 "0.0456"
Set-TextValue 36 5 "  +8.56%  "

# Row 37
Set-TextValue 37 5 "  +0.12%  "

# Row 38
Set-TextValue 38 4 "3.07"
Set-TextValue 38 5 "  +3.02%  "

# Row 39
Set-TextValue 39 4 "17.22"
Set-TextValue 39 5 "  +1.76%  "

# Row 40
Set-TextValue 40 4 "2.61"
Set-TextValue 40 5 "  +1.90%  "

# Row 41
Set-TextValue 41 4 "1.86"
Set-TextValue 41 5 "  +1.81%  "

# Row 42
Set-TextValue 42 5 "  +4.22%  "

# Row 43
Set-TextValue 43 4 "122.62"
Set-TextValue 43 5 "  +3.11%  "

# Row 44
Set-TextValue 44 4 "21.89"
Set-TextValue 44 5 "  +0.58%  "

# Row 45
Set-TextValue 45 4 "0.281"
Set-TextValue 45 5 "  +20.10%  "

# Row 46
Set-TextValue 46 4 "2.05"
Set-TextValue 46 5 "  -1.72%  "

# Row 47
Set-TextValue 47 4 "3.33"
Set-TextValue 47 5 "  +6.50%  "

# Row 49
Set-TextValue 49 4 "2.046.31"
Set-TextValue 49 5 "  +2.00%  "

# Row 50
Set-TextValue 50 5 "  +9.81%  "

# Row 51
Set-TextValue 51 4 "0.874"
Set-TextValue 51 5 "  +3.85%  "
